$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1931216931216931
$ws.Range("C2").Value = 0.5846560846560847
$ws.Range("J2").Value = 0.01851851851851852
$ws.Range("P2").Value = 0.1322751322751323
$ws.Range("S2").Value = 0.07142857142857142
$ws.Range("B3").Value = 0.008583690987124463
$ws.Range("C3").Value = 0.03862660944206009
$ws.Range("J3").Value = 0.008583690987124463
$ws.Range("P3").Value = 0.7725321888412017
$ws.Range("S3").Value = 0.1716738197424893
$ws.Range("B6").Value = 0.05327868852459016
$ws.Range("D6").Value = 0.02049180327868852
$ws.Range("F6").Value = 0.05327868852459016
$ws.Range("J6").Value = 0.2868852459016393
$ws.Range("O6").Value = 0.02459016393442623
$ws.Range("Q6").Value = 0.1680327868852459
$ws.Range("R6").Value = 0.04508196721311476
$ws.Range("S6").Value = 0.3483606557377049
$ws.Range("B7").Value = 0.1141868512110727
$ws.Range("D7").Value = 0.01038062283737024
$ws.Range("F7").Value = 0.04498269896193772
$ws.Range("J7").Value = 0.1349480968858132
$ws.Range("O7").Value = 0.01384083044982699
$ws.Range("Q7").Value = 0.1799307958477509
$ws.Range("R7").Value = 0.05536332179930796
$ws.Range("S7").Value = 0.4463667820069204
$ws.Range("B8").Value = 0.09744779582366589
$ws.Range("D8").Value = 0.02320185614849188
$ws.Range("E8").Value = 0.002320185614849188
$ws.Range("F8").Value = 0.06496519721577726
$ws.Range("J8").Value = 0.1206496519721578
$ws.Range("O8").Value = 0.03248259860788863
$ws.Range("Q8").Value = 0.1740139211136891
$ws.Range("R8").Value = 0.08352668213457076
$ws.Range("S8").Value = 0.4013921113689095
$ws.Range("B9").Value = 0.115
$ws.Range("D9").Value = 0.005
$ws.Range("F9").Value = 0.05
$ws.Range("J9").Value = 0.11
$ws.Range("O9").Value = 0.005
$ws.Range("Q9").Value = 0.185
$ws.Range("R9").Value = 0.08
$ws.Range("S9").Value = 0.45
$ws.Range("B10").Value = 0.1249146757679181
$ws.Range("D10").Value = 0.02116040955631399
$ws.Range("F10").Value = 0.0552901023890785
$ws.Range("J10").Value = 0.1194539249146758
$ws.Range("O10").Value = 0.01979522184300341
$ws.Range("Q10").Value = 0.2177474402730376
$ws.Range("R10").Value = 0.06484641638225255
$ws.Range("S10").Value = 0.3767918088737202
$ws.Range("G11").Value = 0.1574279379157428
$ws.Range("J11").Value = 0.09534368070953436
$ws.Range("K11").Value = 0.1973392461197339
$ws.Range("L11").Value = 0.541019955654102
$ws.Range("S11").Value = 0.008869179600886918
$ws.Range("G12").Value = 0.728
$ws.Range("J12").Value = 0.196
$ws.Range("K12").Value = 0.008
$ws.Range("L12").Value = 0.028
$ws.Range("S12").Value = 0.04
$ws.Range("G13").Value = 0.7704918032786885
$ws.Range("J13").Value = 0.1967213114754098
$ws.Range("S13").Value = 0.03278688524590164
$ws.Range("F15").Value = 0.0310077519379845
$ws.Range("H15").Value = 0.1279069767441861
$ws.Range("I15").Value = 0.07364341085271318
$ws.Range("J15").Value = 0.3488372093023256
$ws.Range("K15").Value = 0.04263565891472868
$ws.Range("M15").Value = 0.02325581395348837
$ws.Range("O15").Value = 0.1085271317829457
$ws.Range("S15").Value = 0.2441860465116279
$ws.Range("F16").Value = 0.03125
$ws.Range("H16").Value = 0.140625
$ws.Range("I16").Value = 0.07421875
$ws.Range("J16").Value = 0.421875
$ws.Range("K16").Value = 0.15234375
$ws.Range("M16").Value = 0.01171875
$ws.Range("O16").Value = 0.04296875
$ws.Range("S16").Value = 0.125
$ws.Range("F17").Value = 0.04069767441860465
$ws.Range("H17").Value = 0.1724806201550388
$ws.Range("I17").Value = 0.07170542635658915
$ws.Range("J17").Value = 0.3992248062015504
$ws.Range("K17").Value = 0.124031007751938
$ws.Range("M17").Value = 0.02131782945736434
$ws.Range("N17").Value = 0.001937984496124031
$ws.Range("O17").Value = 0.05813953488372093
$ws.Range("S17").Value = 0.1104651162790698
$ws.Range("F18").Value = 0.02890173410404624
$ws.Range("H18").Value = 0.1329479768786127
$ws.Range("I18").Value = 0.06358381502890173
$ws.Range("J18").Value = 0.4739884393063584
$ws.Range("K18").Value = 0.1560693641618497
$ws.Range("M18").Value = 0.02890173410404624
$ws.Range("O18").Value = 0.05202312138728324
$ws.Range("S18").Value = 0.06358381502890173
$ws.Range("F19").Value = 0.02534246575342466
$ws.Range("H19").Value = 0.173972602739726
$ws.Range("I19").Value = 0.07876712328767123
$ws.Range("J19").Value = 0.3575342465753424
$ws.Range("K19").Value = 0.1445205479452055
$ws.Range("M19").Value = 0.02671232876712329
$ws.Range("N19").Value = 0.0006849315068493151
$ws.Range("O19").Value = 0.06164383561643835
$ws.Range("S19").Value = 0.1308219178082192
